{"js": "// The goal-of-the-game paragraph:\n//   \"...math question. If...\" -> \"...math questions. If...\"\n//   and \"energy point increase by one\" gets struck through.\nconst questionHits = context.document.body.search(\"math question\", { matchCase: true });\nquestionHits.load(\"text\");\nawait context.sync();\n\nquestionHits.items[0].getRange(\"End\").insertText(\"s\", \"Before\");\nawait context.sync();\n\nconst strikeHits = context.document.body.search(\"energy point increase by one\", { matchCase: true });\nstrikeHits.load(\"text\");\nawait context.sync();\n\nstrikeHits.items[0].font.strikeThrough = true;\nawait context.sync();\n\n// The \"four levels\" paragraph: the text itself does not change, but the\n// \"_GoBack\" bookmark moves from the end of the paragraph to sit right after\n// \"four level\" (before the final \"s\" of \"levels\"), splitting the run there.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst levelHits = context.document.body.search(\"four level\", { matchCase: true });\nlevelHits.load(\"text\");\nawait context.sync();\n\nconst splitPoint = levelHits.items[0].getRange(\"End\").getRange(\"Start\");\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- The goal-of-the-game paragraph ---\n# \"...math question. If...\" -> \"...math questions. If...\"\n$rng = $d.Content\n$rng.Find.Execute(\"math question\") | Out-Null\n$rng.Collapse(0)              # wdCollapseEnd\n$rng.InsertAfter(\"s\")\n\n# Strike-through \"energy point increase by one\"\n$rng2 = $d.Content\n$rng2.Find.Execute(\"energy point increase by one\") | Out-Null\n$rng2.Font.StrikeThrough = 1\n\n# --- The \"four levels\" paragraph ---\n# Text stays the same, but the \"_GoBack\" bookmark moves from the end of the\n# paragraph to right after \"four level\" (before the final \"s\" of \"levels\"),\n# splitting the run at that point.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bm.Delete()\n\n$rng3 = $d.Content\n$rng3.Find.Execute(\"four level\") | Out-Null\n$rng3.Collapse(0)             # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $rng3)\n"}
